$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: replace spaces with underscores in these 4 headers.
$ws.Range("F1").Value = "Fecha_Ingreso"
$ws.Range("I1").Value = "Sueldo_Base"
$ws.Range("J1").Value = "Aportación_ISSSTEESIN"
$ws.Range("K1").Value = "Aportación_Vivienda"

# Update the saved view state: select K18 (clears the prior H1 selection
# and the D1 scroll-freeze offset).
$ws.Range("K18").Select()
